# Auto-generated Excel COM-interop edit script
# Applies updated market-price / profit values to the Raiden_Profits workbook
# as captured by the scheduled runner commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 488.52173
$ws.Range("I33").Value = 340
$ws.Range("K33").Value = 340
$ws.Range("M33").Value = -111
$ws.Range("H40").Value = 5619.933
$ws.Range("J40").Value = 3279.8
$ws.Range("L40").Value = 3279.8
$ws.Range("N40").Value = -3629.8
$ws.Range("H69").Value = 8063.0625
$ws.Range("I69").Value = 6900.9
$ws.Range("K69").Value = 20702.7
$ws.Range("M69").Value = -19828.7
$ws.Range("H72").Value = 8063.0625
$ws.Range("I72").Value = 6900.9
$ws.Range("K72").Value = 62108.1
$ws.Range("M72").Value = -57740.1
$ws.Range("H86").Value = 1300
$ws.Range("I86").Value = 1299
$ws.Range("J86").Value = 1302
$ws.Range("K86").Value = 1299
$ws.Range("L86").Value = 1302
$ws.Range("M86").Value = -176
$ws.Range("N86").Value = -3548
$ws.Range("H89").Value = 1300
$ws.Range("I89").Value = 1299
$ws.Range("J89").Value = 1302
$ws.Range("K89").Value = 6495
$ws.Range("L89").Value = 6510
$ws.Range("M89").Value = -879
$ws.Range("N89").Value = -17742
$ws.Range("H112").Value = 1698.6666
$ws.Range("I112").Value = 1892
$ws.Range("J112").Value = 1689
$ws.Range("K112").Value = 5676
$ws.Range("L112").Value = 5067
$ws.Range("M112").Value = -4568
$ws.Range("N112").Value = -7283
$ws.Range("H132").Value = 2070.2778
$ws.Range("I132").Value = 2250.1428
$ws.Range("J132").Value = 1440.75
$ws.Range("K132").Value = 6750.428400000001
$ws.Range("L132").Value = 4322.25
$ws.Range("M132").Value = -4220.428400000001
$ws.Range("N132").Value = -9382.25
$ws.Range("H137").Value = 3150.6365
$ws.Range("I137").Value = 2557.8
$ws.Range("J137").Value = 3644.6667
$ws.Range("K137").Value = 7673.400000000001
$ws.Range("L137").Value = 10934.0001
$ws.Range("M137").Value = -5123.400000000001
$ws.Range("N137").Value = -16034.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1589.85
$ws.Range("I2").Value = 1540
$ws.Range("J2").Value = 1739.4
$ws.Range("K2").Value = 1540
$ws.Range("L2").Value = 1739.4
$ws.Range("M2").Value = -1427
$ws.Range("N2").Value = -1965.4
$ws.Range("H32").Value = 5909.7417
$ws.Range("I32").Value = 3341.84
$ws.Range("J32").Value = 16609.334
$ws.Range("K32").Value = 3341.84
$ws.Range("L32").Value = 16609.334
$ws.Range("M32").Value = -3054.84
$ws.Range("N32").Value = -17183.334
$ws.Range("H61").Value = 3792.4
$ws.Range("I61").Value = 3365.7778
$ws.Range("J61").Value = 4432.3335
$ws.Range("K61").Value = 3365.7778
$ws.Range("L61").Value = 4432.3335
$ws.Range("M61").Value = -3153.7778
$ws.Range("N61").Value = -4856.3335
$ws.Range("H74").Value = 2020.5416
$ws.Range("J74").Value = 3000.3333
$ws.Range("L74").Value = 3000.3333
$ws.Range("N74").Value = -4748.3333
$ws.Range("H77").Value = 2020.5416
$ws.Range("J77").Value = 3000.3333
$ws.Range("L77").Value = 15001.6665
$ws.Range("N77").Value = -23737.6665
$ws.Range("H102").Value = 2329.3684
$ws.Range("I102").Value = 2474.0667
$ws.Range("J102").Value = 1786.75
$ws.Range("K102").Value = 2474.0667
$ws.Range("L102").Value = 1786.75
$ws.Range("M102").Value = -852.0666999999999
$ws.Range("N102").Value = -5030.75
$ws.Range("H114").Value = 74999.336
$ws.Range("I114").Value = 25000
$ws.Range("J114").Value = 99999
$ws.Range("K114").Value = 25000
$ws.Range("L114").Value = 99999
$ws.Range("M114").Value = -20661
$ws.Range("N114").Value = -108677
$ws.Range("H116").Value = 1589.85
$ws.Range("I116").Value = 1540
$ws.Range("J116").Value = 1739.4
$ws.Range("K116").Value = 1540
$ws.Range("L116").Value = 1739.4
$ws.Range("M116").Value = 754
$ws.Range("N116").Value = -6327.4
$ws.Range("H118").Value = 69408
$ws.Range("J118").Value = 69408
$ws.Range("L118").Value = 69408
$ws.Range("N118").Value = -72722
$ws.Range("H136").Value = 3792.4
$ws.Range("I136").Value = 3365.7778
$ws.Range("J136").Value = 4432.3335
$ws.Range("K136").Value = 10097.3334
$ws.Range("L136").Value = 13297.0005
$ws.Range("M136").Value = -7547.3334
$ws.Range("N136").Value = -18397.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1589.85
$ws.Range("I3").Value = 1540
$ws.Range("J3").Value = 1739.4
$ws.Range("K3").Value = 1540
$ws.Range("L3").Value = 1739.4
$ws.Range("M3").Value = -1426
$ws.Range("N3").Value = -1967.4
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H58").Value = 88421.75
$ws.Range("J58").Value = 92482.57000000001
$ws.Range("L58").Value = 92482.57000000001
$ws.Range("N58").Value = -93070.57000000001
$ws.Range("H60").Value = 61300
$ws.Range("J60").Value = 61300
$ws.Range("L60").Value = 61300
$ws.Range("N60").Value = -62498
$ws.Range("H86").Value = 2979.889
$ws.Range("I86").Value = 2900
$ws.Range("J86").Value = 2989.875
$ws.Range("K86").Value = 2900
$ws.Range("L86").Value = 2989.875
$ws.Range("M86").Value = -1777
$ws.Range("N86").Value = -5235.875
$ws.Range("H89").Value = 2979.889
$ws.Range("I89").Value = 2900
$ws.Range("J89").Value = 2989.875
$ws.Range("K89").Value = 14500
$ws.Range("L89").Value = 14949.375
$ws.Range("M89").Value = -8884
$ws.Range("N89").Value = -26181.375
$ws.Range("H107").Value = 1972.4445
$ws.Range("J107").Value = 2028.909
$ws.Range("L107").Value = 2028.909
$ws.Range("N107").Value = -5868.909
$ws.Range("H134").Value = 3849.9092
$ws.Range("I134").Value = 3131.5386
$ws.Range("J134").Value = 4887.5557
$ws.Range("K134").Value = 9394.6158
$ws.Range("L134").Value = 14662.6671
$ws.Range("M134").Value = -6859.6158
$ws.Range("N134").Value = -19732.6671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3121.125
$ws.Range("I31").Value = 2763.8333
$ws.Range("K31").Value = 2763.8333
$ws.Range("M31").Value = -2468.8333
$ws.Range("H34").Value = 3121.125
$ws.Range("I34").Value = 2763.8333
$ws.Range("K34").Value = 2763.8333
$ws.Range("M34").Value = -2561.8333
$ws.Range("H107").Value = 1064.421
$ws.Range("I107").Value = 757.9091
$ws.Range("K107").Value = 757.9091
$ws.Range("M107").Value = 1162.0909

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 59.5
$ws.Range("I2").Value = 63.454544
$ws.Range("J2").Value = 16
$ws.Range("K2").Value = 380.727264
$ws.Range("L2").Value = 96
$ws.Range("M2").Value = -267.727264
$ws.Range("N2").Value = -322
$ws.Range("H22").Value = 799.5
$ws.Range("I22").Value = 800
$ws.Range("J22").Value = 799
$ws.Range("K22").Value = 2400
$ws.Range("L22").Value = 2397
$ws.Range("M22").Value = -2231
$ws.Range("N22").Value = -2735
$ws.Range("H27").Value = 799.5
$ws.Range("I27").Value = 800
$ws.Range("J27").Value = 799
$ws.Range("K27").Value = 2400
$ws.Range("L27").Value = 2397
$ws.Range("M27").Value = -2298
$ws.Range("N27").Value = -2601
$ws.Range("H58").Value = 7784.3335
$ws.Range("I58").Value = 7799
$ws.Range("J58").Value = 7777
$ws.Range("K58").Value = 23397
$ws.Range("L58").Value = 23331
$ws.Range("M58").Value = -23269
$ws.Range("N58").Value = -23587
$ws.Range("H60").Value = 201.25
$ws.Range("I60").Value = 108
$ws.Range("J60").Value = 214.57143
$ws.Range("K60").Value = 324
$ws.Range("L60").Value = 643.71429
$ws.Range("M60").Value = -73
$ws.Range("N60").Value = -1145.71429
$ws.Range("H113").Value = 617.625
$ws.Range("I113").Value = 398
$ws.Range("J113").Value = 837.25
$ws.Range("K113").Value = 1194
$ws.Range("L113").Value = 2511.75
$ws.Range("M113").Value = 976
$ws.Range("N113").Value = -6851.75
$ws.Range("H122").Value = 713.8333
$ws.Range("I122").Value = 655
$ws.Range("J122").Value = 772.6667
$ws.Range("K122").Value = 5895
$ws.Range("L122").Value = 6954.0003
$ws.Range("M122").Value = -3445
$ws.Range("N122").Value = -11854.0003
$ws.Range("H132").Value = 1745.8334
$ws.Range("I132").Value = 1786.4546
$ws.Range("K132").Value = 16078.0914
$ws.Range("M132").Value = -13548.0914

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7088.4
$ws.Range("J70").Value = 8847.75
$ws.Range("L70").Value = 8847.75
$ws.Range("N70").Value = -9387.75
$ws.Range("H73").Value = 7088.4
$ws.Range("J73").Value = 8847.75
$ws.Range("L73").Value = 8847.75
$ws.Range("N73").Value = -10719.75
$ws.Range("H122").Value = 2900.3125
$ws.Range("I122").Value = 3057
$ws.Range("J122").Value = 2698.8572
$ws.Range("K122").Value = 9171
$ws.Range("L122").Value = 8096.571599999999
$ws.Range("M122").Value = -6721
$ws.Range("N122").Value = -12996.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4527.5
$ws.Range("I7").Value = 4693.0527
$ws.Range("J7").Value = 3898.4
$ws.Range("K7").Value = 4693.0527
$ws.Range("L7").Value = 3898.4
$ws.Range("M7").Value = -4581.0527
$ws.Range("N7").Value = -4122.4
$ws.Range("H46").Value = 3734.8948
$ws.Range("I46").Value = 2997.4546
$ws.Range("J46").Value = 4748.875
$ws.Range("K46").Value = 2997.4546
$ws.Range("L46").Value = 4748.875
$ws.Range("M46").Value = -2809.4546
$ws.Range("N46").Value = -5124.875
$ws.Range("H114").Value = 69999
$ws.Range("J114").Value = 69999
$ws.Range("L114").Value = 69999
$ws.Range("N114").Value = -78677
$ws.Range("H126").Value = 4527.5
$ws.Range("I126").Value = 4693.0527
$ws.Range("J126").Value = 3898.4
$ws.Range("K126").Value = 14079.1581
$ws.Range("L126").Value = 11695.2
$ws.Range("M126").Value = -11609.1581
$ws.Range("N126").Value = -16635.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4073.8667
$ws.Range("I126").Value = 4143.25
$ws.Range("K126").Value = 12429.75
$ws.Range("M126").Value = -9959.75
